# Add a new "CheckBox" worksheet (testcase for checkbox) after "InputForm".
$wb = $excel.ActiveWorkbook

$inputForm = $wb.Worksheets.Item("InputForm")
$ws = $wb.Worksheets.Add($null, $inputForm)
$ws.Name = "CheckBox"

# Row 1 - header
$ws.Range("A1").Value = "checkbox"
$ws.Range("B1").Value = "checkBox"

# Row 2 - checkbox1 (bold, purple, size 9 label)
$ws.Range("A2").Value = "checkbox1"
$ws.Range("B2").Value = "checkBox1"
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").Font.Size = 9
$ws.Range("B2").Font.Color = 7999078

# Row 3 - checkbox2 (bold, purple, size 9 label)
$ws.Range("A3").Value = "checkbox2"
$ws.Range("B3").Value = "checkBox2"
$ws.Range("B3").Font.Bold = $true
$ws.Range("B3").Font.Size = 9
$ws.Range("B3").Font.Color = 7999078

# Row 4 - checkbox3 (plain)
$ws.Range("A4").Value = "checkbox3"
$ws.Range("B4").Value = "checkBox3"

# Row 5 - checkbox4 (bold, purple, size 9 label)
$ws.Range("A5").Value = "checkbox4"
$ws.Range("B5").Value = "checkBox4"
$ws.Range("B5").Font.Bold = $true
$ws.Range("B5").Font.Size = 9
$ws.Range("B5").Font.Color = 7999078

# Row 6 - status / Check All
$ws.Range("A6").Value = "status"
$ws.Range("B6").Value = "Check All"

# Row 7 - Uncheck All
$ws.Range("B7").Value = "Uncheck All"

# Row 9 - success message
$ws.Range("B9").Value = "Success - Check box is checked"

# Select B10 on the new sheet (matches saved selection) and make it the
# active sheet / active tab.
$ws.Range("B10").Select()
$ws.Activate()
